$d = $word.ActiveDocument

# Locate the "ТЕКСТ ДОВОДА" placeholder paragraph (highlighted red text
# followed by a manual line break).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*ТЕКСТ ДОВОДА*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq 0) {
    throw "Could not locate the 'ТЕКСТ ДОВОДА' placeholder paragraph"
}

# 1) Remove the whole "ТЕКСТ ДОВОДА" + line-break paragraph.
$d.Paragraphs.Item($targetIndex).Range.Delete()

# 2) The next paragraph now holds the "{@dovod}" merge-field text; strip the
#    run but keep the (now empty) paragraph and its formatting intact.
$dovodPara = $d.Paragraphs.Item($targetIndex)
$dovodRange = $dovodPara.Range
$d.Range($dovodRange.Start, $dovodRange.End - 1).Delete()

# 3) Insert a brand-new paragraph right after it (inherits the same pPr),
#    which becomes the new home of the "{@dovod}" field and the bookmark.
$dovodPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range
$d.Range($newRange.Start, $newRange.End - 1).Text = "{@dovod}"

# 4) Re-seat the document's "_GoBack" bookmark onto the new paragraph
#    (bookmark names are unique, so re-adding it here removes the old
#    occurrence near the end of the document automatically).
$bmStart = $newRange.Start
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmStart))

# 5) The following (previously-empty) paragraph picks up an explicit
#    en-US language tag on its paragraph mark.
$afterPara = $d.Paragraphs.Item($targetIndex + 2)
$afterPara.Range.LanguageID = "en-US"
